$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Has2"
$ws.Cells.Item(2, 3).Value = "Cd44"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.2004523333333333
$ws.Cells.Item(2, 8).Value = 0.601357
$ws.Cells.Item(2, 9).Value = 0.00350521822469015
$ws.Cells.Item(2, 10).Value = 0.003553987077041381
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 21.084959
$ws.Cells.Item(2, 14).Value = 63.25487699999999
$ws.Cells.Item(2, 15).Value = 0.03381320693734752
$ws.Cells.Item(2, 16).Value = 0.03509122472428063
$ws.Cells.Item(2, 17).Value = 4.226529229787666
$ws.Cells.Item(2, 18).Value = 38.038763068089
$ws.Cells.Item(2, 19).Value = 0.0001185226691920099
$ws.Cells.Item(2, 20).Value = 0.0001247137591876484

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Has2"
$ws.Cells.Item(3, 3).Value = "Cd44"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.2004523333333333
$ws.Cells.Item(3, 8).Value = 0.601357
$ws.Cells.Item(3, 9).Value = 0.00350521822469015
$ws.Cells.Item(3, 10).Value = 0.003553987077041381
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 81.06331633333333
$ws.Cells.Item(3, 14).Value = 243.189949
$ws.Cells.Item(3, 15).Value = 0.12999838843446
$ws.Cells.Item(3, 16).Value = 0.1349118606466557
$ws.Cells.Item(3, 17).Value = 16.24933090675478
$ws.Cells.Item(3, 18).Value = 146.243978160793
$ws.Cells.Item(3, 19).Value = 0.0004556727203208184
$ws.Cells.Item(3, 20).Value = 0.0004794750092778219

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Has2"
$ws.Cells.Item(4, 3).Value = "Cd44"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.2004523333333333
$ws.Cells.Item(4, 8).Value = 0.601357
$ws.Cells.Item(4, 9).Value = 0.00350521822469015
$ws.Cells.Item(4, 10).Value = 0.003553987077041381
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 136.9994176666667
$ws.Cells.Item(4, 14).Value = 410.998253
$ws.Cells.Item(4, 15).Value = 0.2197011461990087
$ws.Cells.Item(4, 16).Value = 0.2280050605000741
$ws.Cells.Item(4, 17).Value = 27.46185293659122
$ws.Cells.Item(4, 18).Value = 247.156676429321
$ws.Cells.Item(4, 19).Value = 0.0007701004616420804
$ws.Cells.Item(4, 20).Value = 0.0008103270385173016

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Has2"
$ws.Cells.Item(5, 3).Value = "Cd44"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.2004523333333333
$ws.Cells.Item(5, 8).Value = 0.601357
$ws.Cells.Item(5, 9).Value = 0.00350521822469015
$ws.Cells.Item(5, 10).Value = 0.003553987077041381
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 316.292811
$ws.Cells.Item(5, 14).Value = 948.878433
$ws.Cells.Item(5, 15).Value = 0.5072276531881493
$ws.Cells.Item(5, 16).Value = 0.5263990368430604
$ws.Cells.Item(5, 17).Value = 63.401631981509
$ws.Cells.Item(5, 18).Value = 570.614687833581
$ws.Cells.Item(5, 19).Value = 0.001777943614021916
$ws.Cells.Item(5, 20).Value = 0.001870815374307267

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Has2"
$ws.Cells.Item(6, 3).Value = "Cd44"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.2004523333333333
$ws.Cells.Item(6, 8).Value = 0.601357
$ws.Cells.Item(6, 9).Value = 0.00350521822469015
$ws.Cells.Item(6, 10).Value = 0.003553987077041381
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 68.131198
$ws.Cells.Item(6, 14).Value = 136.262396
$ws.Cells.Item(6, 15).Value = 0.1092596052410345
$ws.Cells.Item(6, 16).Value = 0.07559281728592908
$ws.Cells.Item(6, 17).Value = 13.65705761189533
$ws.Cells.Item(6, 18).Value = 81.942345671372
$ws.Cells.Item(6, 19).Value = 0.0003829787595133256
$ws.Cells.Item(6, 20).Value = 0.0002686558957513423

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Has2"
$ws.Cells.Item(7, 3).Value = "Cd44"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 54.63217833333334
$ws.Cells.Item(7, 8).Value = 163.896535
$ws.Cells.Item(7, 9).Value = 0.9553279024698591
$ws.Cells.Item(7, 10).Value = 0.9686195843099198
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 21.084959
$ws.Cells.Item(7, 14).Value = 63.25487699999999
$ws.Cells.Item(7, 15).Value = 0.03381320693734752
$ws.Cells.Item(7, 16).Value = 0.03509122472428063
$ws.Cells.Item(7, 17).Value = 1151.917240239022
$ws.Cells.Item(7, 18).Value = 10367.25516215119
$ws.Cells.Item(7, 19).Value = 0.0323027000592355
$ws.Cells.Item(7, 20).Value = 0.03399004750535868

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Has2"
$ws.Cells.Item(8, 3).Value = "Cd44"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 54.63217833333334
$ws.Cells.Item(8, 8).Value = 163.896535
$ws.Cells.Item(8, 9).Value = 0.9553279024698591
$ws.Cells.Item(8, 10).Value = 0.9686195843099198
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 81.06331633333333
$ws.Cells.Item(8, 14).Value = 243.189949
$ws.Cells.Item(8, 15).Value = 0.12999838843446
$ws.Cells.Item(8, 16).Value = 0.1349118606466557
$ws.Cells.Item(8, 17).Value = 4428.66555421408
$ws.Cells.Item(8, 18).Value = 39857.98998792672
$ws.Cells.Item(8, 19).Value = 0.1241910877475546
$ws.Cells.Item(8, 20).Value = 0.1306782703780414

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Has2"
$ws.Cells.Item(9, 3).Value = "Cd44"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 54.63217833333334
$ws.Cells.Item(9, 8).Value = 163.896535
$ws.Cells.Item(9, 9).Value = 0.9553279024698591
$ws.Cells.Item(9, 10).Value = 0.9686195843099198
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 136.9994176666667
$ws.Cells.Item(9, 14).Value = 410.998253
$ws.Cells.Item(9, 15).Value = 0.2197011461990087
$ws.Cells.Item(9, 16).Value = 0.2280050605000741
$ws.Cells.Item(9, 17).Value = 7484.576617528151
$ws.Cells.Item(9, 18).Value = 67361.18955775336
$ws.Cells.Item(9, 19).Value = 0.2098866351685228
$ws.Cells.Item(9, 20).Value = 0.2208501669221399

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Has2"
$ws.Cells.Item(10, 3).Value = "Cd44"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 54.63217833333334
$ws.Cells.Item(10, 8).Value = 163.896535
$ws.Cells.Item(10, 9).Value = 0.9553279024698591
$ws.Cells.Item(10, 10).Value = 0.9686195843099198
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 316.292811
$ws.Cells.Item(10, 14).Value = 948.878433
$ws.Cells.Item(10, 15).Value = 0.5072276531881493
$ws.Cells.Item(10, 16).Value = 0.5263990368430604
$ws.Cells.Item(10, 17).Value = 17279.7652561033
$ws.Cells.Item(10, 18).Value = 155517.8873049296
$ws.Cells.Item(10, 19).Value = 0.4845687299949438
$ws.Cells.Item(10, 20).Value = 0.5098804162480673

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Has2"
$ws.Cells.Item(11, 3).Value = "Cd44"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 54.63217833333334
$ws.Cells.Item(11, 8).Value = 163.896535
$ws.Cells.Item(11, 9).Value = 0.9553279024698591
$ws.Cells.Item(11, 10).Value = 0.9686195843099198
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 68.131198
$ws.Cells.Item(11, 14).Value = 136.262396
$ws.Cells.Item(11, 15).Value = 0.1092596052410345
$ws.Cells.Item(11, 16).Value = 0.07559281728592908
$ws.Cells.Item(11, 17).Value = 3722.155759199643
$ws.Cells.Item(11, 18).Value = 22332.93455519786
$ws.Cells.Item(11, 19).Value = 0.1043787494996023
$ws.Cells.Item(11, 20).Value = 0.07322068325631234

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Has2"
$ws.Cells.Item(12, 3).Value = "Cd44"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 2.3542035
$ws.Cells.Item(12, 8).Value = 4.708407
$ws.Cells.Item(12, 9).Value = 0.04116687930545086
$ws.Cells.Item(12, 10).Value = 0.02782642861303881
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 21.084959
$ws.Cells.Item(12, 14).Value = 63.25487699999999
$ws.Cells.Item(12, 15).Value = 0.03381320693734752
$ws.Cells.Item(12, 16).Value = 0.03509122472428063
$ws.Cells.Item(12, 17).Value = 49.6382842751565
$ws.Cells.Item(12, 18).Value = 297.829705650939
$ws.Cells.Item(12, 19).Value = 0.001391984208920019
$ws.Cells.Item(12, 20).Value = 0.0009764634597342974

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Has2"
$ws.Cells.Item(13, 3).Value = "Cd44"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.3542035
$ws.Cells.Item(13, 8).Value = 4.708407
$ws.Cells.Item(13, 9).Value = 0.04116687930545086
$ws.Cells.Item(13, 10).Value = 0.02782642861303881
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 81.06331633333333
$ws.Cells.Item(13, 14).Value = 243.189949
$ws.Cells.Item(13, 15).Value = 0.12999838843446
$ws.Cells.Item(13, 16).Value = 0.1349118606466557
$ws.Cells.Item(13, 17).Value = 190.8395430335405
$ws.Cells.Item(13, 18).Value = 1145.037258201243
$ws.Cells.Item(13, 19).Value = 0.005351627966584534
$ws.Cells.Item(13, 20).Value = 0.003754115259336404

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Has2"
$ws.Cells.Item(14, 3).Value = "Cd44"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.3542035
$ws.Cells.Item(14, 8).Value = 4.708407
$ws.Cells.Item(14, 9).Value = 0.04116687930545086
$ws.Cells.Item(14, 10).Value = 0.02782642861303881
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 136.9994176666667
$ws.Cells.Item(14, 14).Value = 410.998253
$ws.Cells.Item(14, 15).Value = 0.2197011461990087
$ws.Cells.Item(14, 16).Value = 0.2280050605000741
$ws.Cells.Item(14, 17).Value = 322.5245085688285
$ws.Cells.Item(14, 18).Value = 1935.147051412971
$ws.Cells.Item(14, 19).Value = 0.009044410568843805
$ws.Cells.Item(14, 20).Value = 0.006344566539416906

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Has2"
$ws.Cells.Item(15, 3).Value = "Cd44"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.3542035
$ws.Cells.Item(15, 8).Value = 4.708407
$ws.Cells.Item(15, 9).Value = 0.04116687930545086
$ws.Cells.Item(15, 10).Value = 0.02782642861303881
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 316.292811
$ws.Cells.Item(15, 14).Value = 948.878433
$ws.Cells.Item(15, 15).Value = 0.5072276531881493
$ws.Cells.Item(15, 16).Value = 0.5263990368430604
$ws.Cells.Item(15, 17).Value = 744.6176426810384
$ws.Cells.Item(15, 18).Value = 4467.705856086232
$ws.Cells.Item(15, 19).Value = 0.02088097957918363
$ws.Cells.Item(15, 20).Value = 0.01464780522068581

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Has2"
$ws.Cells.Item(16, 3).Value = "Cd44"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.3542035
$ws.Cells.Item(16, 8).Value = 4.708407
$ws.Cells.Item(16, 9).Value = 0.04116687930545086
$ws.Cells.Item(16, 10).Value = 0.02782642861303881
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 68.131198
$ws.Cells.Item(16, 14).Value = 136.262396
$ws.Cells.Item(16, 15).Value = 0.1092596052410345
$ws.Cells.Item(16, 16).Value = 0.07559281728592908
$ws.Cells.Item(16, 17).Value = 160.394704790793
$ws.Cells.Item(16, 18).Value = 641.578819163172
$ws.Cells.Item(16, 19).Value = 0.004497876981918873
$ws.Cells.Item(16, 20).Value = 0.002103478133865392
